$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = 800
$ws.Range("B15").Value = 840
$ws.Range("B16").Value = 4900
$ws.Range("B19").Value = 70
$ws.Range("B20").Value = 260
